# Add four more lab entries below the existing "Name" header on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Lab 9"
$ws.Range("A3").Value = "Lab 12"
$ws.Range("A4").Value = "Lab 13"
$ws.Range("A5").Value = "Lab 14"

# Move the active selection below the newly written data, as in the source.
$ws.Range("A6").Select() | Out-Null
